$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 1184.4445
$ws.Range("I40").Value = 1030
$ws.Range("J40").Value = 1409.091
$ws.Range("K40").Value = 1030
$ws.Range("L40").Value = 1409.091
$ws.Range("M40").Value = -855
$ws.Range("N40").Value = -1759.091

# Row 92
$ws.Range("H92").Value = 794.1579
$ws.Range("I92").Value = 794.1579
$ws.Range("K92").Value = 794.1579
$ws.Range("M92").Value = 453.8421

# Row 133
$ws.Range("H133").Value = 43599.4
$ws.Range("J133").Value = 43599.4
$ws.Range("L133").Value = 43599.4
$ws.Range("N133").Value = -53719.4

# Row 137
$ws.Range("H137").Value = 1151.091
$ws.Range("I137").Value = 1045.7894
$ws.Range("J137").Value = 1294
$ws.Range("K137").Value = 3137.3682
$ws.Range("L137").Value = 3882
$ws.Range("M137").Value = -587.3681999999999
$ws.Range("N137").Value = -8982

# Row 138
$ws.Range("H138").Value = 2264.17
$ws.Range("I138").Value = 1178
$ws.Range("J138").Value = 2321.337
$ws.Range("K138").Value = 3534
$ws.Range("L138").Value = 6964.011
$ws.Range("M138").Value = 1606
$ws.Range("N138").Value = -17244.011

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 24
$ws.Range("H24").Value = 13999.5
$ws.Range("J24").Value = 13999.5
$ws.Range("L24").Value = 13999.5
$ws.Range("N24").Value = -14747.5

# Row 32
$ws.Range("H32").Value = 5953.5454
$ws.Range("I32").Value = 6189.476
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 6189.476
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -5902.476
$ws.Range("N32").Value = -1573

# Row 54
$ws.Range("N54").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0

# Row 61
$ws.Range("H61").Value = 1578.25
$ws.Range("I61").Value = 1449
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1449
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1237
$ws.Range("N61").Value = -3424

# Row 74
$ws.Range("H74").Value = 798.1111
$ws.Range("I74").Value = 781.3913
$ws.Range("J74").Value = 894.25
$ws.Range("K74").Value = 781.3913
$ws.Range("L74").Value = 894.25
$ws.Range("M74").Value = 92.6087
$ws.Range("N74").Value = -2642.25

# Row 77
$ws.Range("H77").Value = 798.1111
$ws.Range("I77").Value = 781.3913
$ws.Range("J77").Value = 894.25
$ws.Range("K77").Value = 3906.9565
$ws.Range("L77").Value = 4471.25
$ws.Range("M77").Value = 461.0434999999998
$ws.Range("N77").Value = -13207.25

# Row 100
$ws.Range("H100").Value = 13999.5
$ws.Range("J100").Value = 13999.5
$ws.Range("L100").Value = 13999.5
$ws.Range("N100").Value = -16163.5

# Row 102
$ws.Range("H102").Value = 16677658
$ws.Range("I102").Value = 20846598
$ws.Range("J102").Value = 1900
$ws.Range("K102").Value = 20846598
$ws.Range("L102").Value = 1900
$ws.Range("M102").Value = -20844976
$ws.Range("N102").Value = -5144

# Row 110
$ws.Range("H110").Value = 999.7857
$ws.Range("I110").Value = 768.5
$ws.Range("J110").Value = 4006.5
$ws.Range("K110").Value = 768.5
$ws.Range("L110").Value = 4006.5
$ws.Range("M110").Value = 1276.5
$ws.Range("N110").Value = -8096.5

# Row 122
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 1291.0588
$ws.Range("I122").Value = 1291.0588
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3873.1764
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1423.1764

# Row 136
$ws.Range("H136").Value = 1578.25
$ws.Range("I136").Value = 1449
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4347
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1797
$ws.Range("N136").Value = -14100

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 3689.3157
$ws.Range("I86").Value = 3665.7058
$ws.Range("J86").Value = 3890
$ws.Range("K86").Value = 3665.7058
$ws.Range("L86").Value = 3890
$ws.Range("M86").Value = -2542.7058
$ws.Range("N86").Value = -6136

# Row 89
$ws.Range("H89").Value = 3689.3157
$ws.Range("I89").Value = 3665.7058
$ws.Range("J89").Value = 3890
$ws.Range("K89").Value = 18328.529
$ws.Range("L89").Value = 19450
$ws.Range("M89").Value = -12712.529
$ws.Range("N89").Value = -30682

# Row 105
$ws.Range("H105").Value = 83335720
$ws.Range("I105").Value = 142859360
$ws.Range("J105").Value = 2620
$ws.Range("K105").Value = 142859360
$ws.Range("L105").Value = 2620
$ws.Range("M105").Value = -142857613
$ws.Range("N105").Value = -6114

# Row 126
$ws.Range("H126").Value = 50000
$ws.Range("J126").Value = 50000
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880

# Row 132
$ws.Range("H132").Value = 30779
$ws.Range("J132").Value = 30779
$ws.Range("L132").Value = 30779
$ws.Range("N132").Value = -40899

# Row 134
$ws.Range("H134").Value = 9319.923000000001
$ws.Range("I134").Value = 6425.95
$ws.Range("J134").Value = 18966.5
$ws.Range("K134").Value = 19277.85
$ws.Range("L134").Value = 56899.5
$ws.Range("M134").Value = -16742.85
$ws.Range("N134").Value = -61969.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 763.1622
$ws.Range("I31").Value = 672.5848999999999
$ws.Range("J31").Value = 991.7619
$ws.Range("K31").Value = 672.5848999999999
$ws.Range("L31").Value = 991.7619
$ws.Range("M31").Value = -377.5848999999999
$ws.Range("N31").Value = -1581.7619

# Row 34
$ws.Range("H34").Value = 763.1622
$ws.Range("I34").Value = 672.5848999999999
$ws.Range("J34").Value = 991.7619
$ws.Range("K34").Value = 672.5848999999999
$ws.Range("L34").Value = 991.7619
$ws.Range("M34").Value = -470.5848999999999
$ws.Range("N34").Value = -1395.7619

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 68
$ws.Range("H68").Value = 1436.3334
$ws.Range("I68").Value = 683.5833
$ws.Range("J68").Value = 2038.5333
$ws.Range("K68").Value = 2050.7499
$ws.Range("L68").Value = 6115.5999
$ws.Range("M68").Value = -1239.7499
$ws.Range("N68").Value = -7737.5999

# Row 69
$ws.Range("M69").ClearContents()
$ws.Range("H69").Value = 3402
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3402
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 10206
$ws.Range("N69").Value = -11828

# Row 71
$ws.Range("H71").Value = 1436.3334
$ws.Range("I71").Value = 683.5833
$ws.Range("J71").Value = 2038.5333
$ws.Range("K71").Value = 6152.2497
$ws.Range("L71").Value = 18346.7997
$ws.Range("M71").Value = -2096.2497
$ws.Range("N71").Value = -26458.7997

# Row 72
$ws.Range("M72").ClearContents()
$ws.Range("H72").Value = 3402
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3402
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 30618
$ws.Range("N72").Value = -38730

# Row 134
$ws.Range("H134").Value = 4447.8667
$ws.Range("I134").Value = 2302.375
$ws.Range("J134").Value = 6899.857
$ws.Range("K134").Value = 6907.125
$ws.Range("L134").Value = 20699.571
$ws.Range("M134").Value = -1837.125
$ws.Range("N134").Value = -30839.571

# Row 141
$ws.Range("H141").Value = 2513.3333
$ws.Range("I141").Value = 2327.5
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 6982.5
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -1802.5
$ws.Range("N141").Value = -22360

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 52
$ws.Range("H52").Value = 10515
$ws.Range("I52").Value = 1030
$ws.Range("J52").Value = 20000
$ws.Range("K52").Value = 1030
$ws.Range("L52").Value = 20000
$ws.Range("M52").Value = -771
$ws.Range("N52").Value = -20518

# Row 126
$ws.Range("H126").Value = 2636.3333
$ws.Range("I126").Value = 1806
$ws.Range("J126").Value = 3466.6667
$ws.Range("K126").Value = 5418
$ws.Range("L126").Value = 10400.0001
$ws.Range("M126").Value = -2948
$ws.Range("N126").Value = -15340.0001

# Row 132
$ws.Range("H132").Value = 2191.111
$ws.Range("I132").Value = 1812.2609
$ws.Range("J132").Value = 4369.5
$ws.Range("K132").Value = 5436.7827
$ws.Range("L132").Value = 13108.5
$ws.Range("M132").Value = -2906.7827
$ws.Range("N132").Value = -18168.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 35
$ws.Range("H35").Value = 3524.5
$ws.Range("I35").Value = 3524.5
$ws.Range("K35").Value = 3524.5
$ws.Range("M35").Value = -3188.5

# Row 40
$ws.Range("H40").Value = 2398.625
$ws.Range("I40").Value = 2198.3333
$ws.Range("K40").Value = 2198.3333
$ws.Range("M40").Value = -2062.3333

# Row 45
$ws.Range("H45").Value = 3000
$ws.Range("I45").Value = 3000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2593

# Row 46
$ws.Range("H46").Value = 1433.3334
$ws.Range("J46").Value = 1900
$ws.Range("L46").Value = 1900
$ws.Range("N46").Value = -2276

# Row 68
$ws.Range("H68").Value = 1463.44
$ws.Range("I68").Value = 1241.5264
$ws.Range("J68").Value = 2166.1667
$ws.Range("K68").Value = 1241.5264
$ws.Range("L68").Value = 2166.1667
$ws.Range("M68").Value = -492.5264
$ws.Range("N68").Value = -3664.1667

# Row 71
$ws.Range("H71").Value = 1463.44
$ws.Range("I71").Value = 1241.5264
$ws.Range("J71").Value = 2166.1667
$ws.Range("K71").Value = 6207.632
$ws.Range("L71").Value = 10830.8335
$ws.Range("M71").Value = -2463.632
$ws.Range("N71").Value = -18318.8335

# Row 122
$ws.Range("H122").Value = 35716660
$ws.Range("I122").Value = 50002120
$ws.Range("J122").Value = 3002.5
$ws.Range("K122").Value = 150006360
$ws.Range("L122").Value = 9007.5
$ws.Range("M122").Value = -150003910
$ws.Range("N122").Value = -13907.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 52
$ws.Range("M52").ClearContents()
$ws.Range("H52").Value = 24589
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 24589
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 24589
$ws.Range("N52").Value = -25041

# Row 62
$ws.Range("H62").Value = 62505476
$ws.Range("I62").Value = 71432330
$ws.Range("J62").Value = 17500
$ws.Range("K62").Value = 71432330
$ws.Range("L62").Value = 17500
$ws.Range("M62").Value = -71431706
$ws.Range("N62").Value = -18748

# Row 65
$ws.Range("H65").Value = 62505476
$ws.Range("I65").Value = 71432330
$ws.Range("J65").Value = 17500
$ws.Range("K65").Value = 357161650
$ws.Range("L65").Value = 87500
$ws.Range("M65").Value = -357158530
$ws.Range("N65").Value = -93740

# Row 122
$ws.Range("H122").Value = 9631195
$ws.Range("I122").Value = 11819712
$ws.Range("K122").Value = 35459136
$ws.Range("M122").Value = -35456686
